$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1041
$ws.Range("F6").Value = 371
$ws.Range("F13").Value = 523
$ws.Range("F18").Value = 1427
$ws.Range("F21").Value = 1159
$ws.Range("F22").Value = 36
$ws.Range("F23").Value = 418
$ws.Range("F24").Value = 38
$ws.Range("F25").Value = 3597
$ws.Range("F26").Value = 712
$ws.Range("F27").Value = 566
$ws.Range("F28").Value = 1589
$ws.Range("F29").Value = 46

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 17
$ws.Range("F8").Value = 34
$ws.Range("F9").Value = 29
$ws.Range("F13").Value = 89

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 23

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 23
$ws.Range("F9").Value = 17
$ws.Range("F13").Value = 34
$ws.Range("F14").Value = 29
$ws.Range("F16").Value = 1041
$ws.Range("F17").Value = 371
$ws.Range("F24").Value = 523
$ws.Range("F29").Value = 1427
$ws.Range("F34").Value = 1159
$ws.Range("F35").Value = 36
$ws.Range("F36").Value = 418
$ws.Range("F37").Value = 38
$ws.Range("F38").Value = 3597
$ws.Range("F39").Value = 712
$ws.Range("F40").Value = 566
$ws.Range("F41").Value = 1589
$ws.Range("F42").Value = 89
$ws.Range("F44").Value = 46
